$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 12500
$ws.Range("I2").Value = 10000
$ws.Range("S2").Value = 33061
$ws.Range("T2").Value = 5005.48337366473
$ws.Range("U2").Value = 42500
$ws.Range("X2").Value = 0

# Row 3
$ws.Range("P3").Value = 0
$ws.Range("S3").Value = 31828
$ws.Range("T3").Value = 4922.298456099087
$ws.Range("U3").Value = 22500

# Row 4
$ws.Range("S4").Value = 30574
$ws.Range("T4").Value = 4850.625419333438

# Row 5
$ws.Range("S5").Value = 29643
$ws.Range("T5").Value = 4797.274898040722

# Row 6
$ws.Range("P6").Value = 0
$ws.Range("S6").Value = 29848
$ws.Range("T6").Value = 4859.095012476347
$ws.Range("U6").Value = 22500

# Row 7
$ws.Range("P7").Value = 0
$ws.Range("S7").Value = 30796
$ws.Range("T7").Value = 4960.889853826959
$ws.Range("U7").Value = 22500

# Row 8
$ws.Range("P8").Value = 0
$ws.Range("S8").Value = 31613
$ws.Range("T8").Value = 5211.449593293681
$ws.Range("U8").Value = 22500

# Row 9
$ws.Range("P9").Value = 0
$ws.Range("S9").Value = 34939
$ws.Range("T9").Value = 6049.649550323782
$ws.Range("U9").Value = 22500

# Row 10
$ws.Range("S10").Value = 41934
$ws.Range("T10").Value = 8417.127547097272

# Row 11
$ws.Range("S11").Value = 43152
$ws.Range("T11").Value = 13902.175

# Row 12
$ws.Range("S12").Value = 44379
$ws.Range("T12").Value = 15210.828

# Row 13
$ws.Range("S13").Value = 45163
$ws.Range("T13").Value = 15097.205

# Row 14
$ws.Range("S14").Value = 45347
$ws.Range("T14").Value = 14468.754

# Row 15
$ws.Range("S15").Value = 47292
$ws.Range("T15").Value = 14615.9345
$ws.Range("V15").Value = 2988.007433333334
$ws.Range("W15").Value = 8.447192755792639

# Row 16 (S16 removed entirely)
$ws.Range("S16").Value = ""
$ws.Range("T16").Value = 15642.9465

# Row 17 (S17 removed entirely)
$ws.Range("S17").Value = ""
$ws.Range("T17").Value = 15864.0115

# Row 18 (S18 removed entirely)
$ws.Range("S18").Value = ""
$ws.Range("T18").Value = 16275.092

# Row 19
$ws.Range("T19").Value = 16265.2175

# Row 20
$ws.Range("T20").Value = 15390.046

# Row 21
$ws.Range("T21").Value = 13874.261

# Row 22
$ws.Range("T22").Value = 11951.5345

# Row 23
$ws.Range("B23").Value = 24000
$ws.Range("T23").Value = 9041.0965
$ws.Range("U23").Value = 74000

# Row 24
$ws.Range("B24").Value = 21000
$ws.Range("T24").Value = 6150.5125
$ws.Range("U24").Value = 71000

# Row 25
$ws.Range("B25").Value = 17500
$ws.Range("T25").Value = 5731.305
$ws.Range("U25").Value = 67500
